$d = $word.ActiveDocument
$t = $d.Tables(1)

# Map of 1-based table row index -> new cell text.
# Rows 1-12 get their single numeric value replaced.
# Rows 44-46 previously held a full tab-separated line of stats; they are
# collapsed back down to the single leading value (the rest of that line
# moved/duplicated data already present higher up the table).
$rowValues = @{
    1  = "0M"
    2  = "0M"
    3  = "0M"
    4  = "2000"
    5  = "0.00003"
    6  = "0.54812"
    7  = "0.07619"
    8  = "0.02168"
    9  = "0.36450"
    10 = "0.36450"
    11 = "0.54812"
    12 = "1.28801"
    44 = "99.41"
    45 = "1.29"
    46 = "218"
}

foreach ($rowIndex in $rowValues.Keys) {
    $cell = $t.Rows($rowIndex).Cells(1)
    $cell.Range.Text = $rowValues[$rowIndex]
}
